$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003455333333333333
$ws.Range("H2").Value = 0.010366
$ws.Range("I2").Value = 0.000270121469710956
$ws.Range("J2").Value = 0.000270121469710956
$ws.Range("M2").Value = 0.092903
$ws.Range("N2").Value = 0.278709
$ws.Range("O2").Value = 0.03600043090620505
$ws.Range("P2").Value = 0.03600043090620505
$ws.Range("Q2").Value = 0.0003210108326666667
$ws.Range("R2").Value = 0.002889097494
$ws.Range("S2").Value = 0.00000972448930661183
$ws.Range("T2").Value = 0.000009724489306611832
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003455333333333333
$ws.Range("H3").Value = 0.010366
$ws.Range("I3").Value = 0.000270121469710956
$ws.Range("J3").Value = 0.000270121469710956
$ws.Range("O3").Value = 0.1132051051535142
$ws.Range("P3").Value = 0.1132051051535142
$ws.Range("Q3").Value = 0.001009434169333333
$ws.Range("R3").Value = 0.009084907524000001
$ws.Range("S3").Value = 0.00003057912938285057
$ws.Range("T3").Value = 0.00003057912938285058
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.003455333333333333
$ws.Range("H4").Value = 0.010366
$ws.Range("I4").Value = 0.000270121469710956
$ws.Range("J4").Value = 0.000270121469710956
$ws.Range("M4").Value = 2.195567
$ws.Range("N4").Value = 6.586701
$ws.Range("O4").Value = 0.8507944639402807
$ws.Range("P4").Value = 0.8507944639402808
$ws.Range("Q4").Value = 0.007586415840666666
$ws.Range("R4").Value = 0.068277742566
$ws.Range("S4").Value = 0.0002298178510214935
$ws.Range("T4").Value = 0.0002298178510214936
$ws.Range("I5").Value = 0.8991926531546518
$ws.Range("J5").Value = 0.8991926531546519
$ws.Range("M5").Value = 0.092903
$ws.Range("N5").Value = 0.278709
$ws.Range("O5").Value = 0.03600043090620505
$ws.Range("P5").Value = 0.03600043090620505
$ws.Range("Q5").Value = 1.068595482713
$ws.Range("R5").Value = 9.617359344416998
$ws.Range("S5").Value = 0.03237132298126124
$ws.Range("T5").Value = 0.03237132298126125
$ws.Range("I6").Value = 0.8991926531546518
$ws.Range("J6").Value = 0.8991926531546519
$ws.Range("O6").Value = 0.1132051051535142
$ws.Range("P6").Value = 0.1132051051535142
$ws.Range("S6").Value = 0.1017931988536398
$ws.Range("T6").Value = 0.1017931988536398
$ws.Range("I7").Value = 0.8991926531546518
$ws.Range("J7").Value = 0.8991926531546519
$ws.Range("M7").Value = 2.195567
$ws.Range("N7").Value = 6.586701
$ws.Range("O7").Value = 0.8507944639402807
$ws.Range("P7").Value = 0.8507944639402808
$ws.Range("Q7").Value = 25.254006632657
$ws.Range("R7").Value = 227.286059693913
$ws.Range("S7").Value = 0.7650281313197507
$ws.Range("T7").Value = 0.7650281313197509
$ws.Range("G8").Value = 1.286049666666667
$ws.Range("H8").Value = 3.858149
$ws.Range("I8").Value = 0.1005372253756372
$ws.Range("J8").Value = 0.1005372253756372
$ws.Range("M8").Value = 0.092903
$ws.Range("N8").Value = 0.278709
$ws.Range("O8").Value = 0.03600043090620505
$ws.Range("P8").Value = 0.03600043090620505
$ws.Range("Q8").Value = 0.1194778721823333
$ws.Range("R8").Value = 1.075300849641
$ws.Range("S8").Value = 0.003619383435637192
$ws.Range("T8").Value = 0.003619383435637192
$ws.Range("G9").Value = 1.286049666666667
$ws.Range("H9").Value = 3.858149
$ws.Range("I9").Value = 0.1005372253756372
$ws.Range("J9").Value = 0.1005372253756372
$ws.Range("O9").Value = 0.1132051051535142
$ws.Range("P9").Value = 0.1132051051535142
$ws.Range("Q9").Value = 0.3757039775206667
$ws.Range("R9").Value = 3.381335797686
$ws.Range("S9").Value = 0.01138132717049157
$ws.Range("T9").Value = 0.01138132717049157
$ws.Range("G10").Value = 1.286049666666667
$ws.Range("H10").Value = 3.858149
$ws.Range("I10").Value = 0.1005372253756372
$ws.Range("J10").Value = 0.1005372253756372
$ws.Range("M10").Value = 2.195567
$ws.Range("N10").Value = 6.586701
$ws.Range("O10").Value = 0.8507944639402807
$ws.Range("P10").Value = 0.8507944639402808
$ws.Range("Q10").Value = 2.823608208494333
$ws.Range("R10").Value = 25.412473876449
$ws.Range("S10").Value = 0.08553651476950844
$ws.Range("T10").Value = 0.08553651476950845
